# Add the new "2022-Q3" sheet (inserted right after "总计", before "2022-Q2")
# and update the "总计" summary sheet with a new leading row for 2022-Q3.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet positioned before "2022-Q2"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($wsQ2)
$ws.Name = "2022-Q3"

# Copy header formatting (bold + border + centered) from the 总计 sheet header
$wsTotal.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Copy the "index column" number formatting (bold + border + centered) for column A
$wsTotal.Range("A2").Copy()
$ws.Range("A2:A18").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Force columns B:G (rows 2-18) to be stored as text so values such as
# "009556" or "63.72" keep their original literal representation.
$ws.Range("B2:G18").NumberFormat = "@"

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "009556"
$ws.Cells.Item(2,3).Value = "兴全合丰三年持有期混合"
$ws.Cells.Item(2,4).Value = "63.72"
$ws.Cells.Item(2,5).Value = "89.42"
$ws.Cells.Item(2,6).Value = "3.97"
$ws.Cells.Item(2,7).Value = "2.5297"
$ws.Cells.Item(2,8).Value = 8

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "015188"
$ws.Cells.Item(3,3).Value = "汇添富消费升级混合D"
$ws.Cells.Item(3,4).Value = "48.66"
$ws.Cells.Item(3,5).Value = "87.74"
$ws.Cells.Item(3,6).Value = "2.36"
$ws.Cells.Item(3,7).Value = "1.1484"
$ws.Cells.Item(3,8).Value = 10

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "006408"
$ws.Cells.Item(4,3).Value = "汇添富消费升级混合A"
$ws.Cells.Item(4,4).Value = "48.10"
$ws.Cells.Item(4,5).Value = "87.74"
$ws.Cells.Item(4,6).Value = "2.36"
$ws.Cells.Item(4,7).Value = "1.1352"
$ws.Cells.Item(4,8).Value = 10

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "005644"
$ws.Cells.Item(5,3).Value = "广发沪港深行业龙头混合"
$ws.Cells.Item(5,4).Value = "12.35"
$ws.Cells.Item(5,5).Value = "85.09"
$ws.Cells.Item(5,6).Value = "5.16"
$ws.Cells.Item(5,7).Value = "0.6373"
$ws.Cells.Item(5,8).Value = 6

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "004424"
$ws.Cells.Item(6,3).Value = "汇添富文体娱乐主题混合A"
$ws.Cells.Item(6,4).Value = "19.35"
$ws.Cells.Item(6,5).Value = "90.70"
$ws.Cells.Item(6,6).Value = "2.92"
$ws.Cells.Item(6,7).Value = "0.5650"
$ws.Cells.Item(6,8).Value = 10

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "006595"
$ws.Cells.Item(7,3).Value = "广发港股通优质增长混合A"
$ws.Cells.Item(7,4).Value = "7.56"
$ws.Cells.Item(7,5).Value = "90.07"
$ws.Cells.Item(7,6).Value = "5.56"
$ws.Cells.Item(7,7).Value = "0.4203"
$ws.Cells.Item(7,8).Value = 10

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "012640"
$ws.Cells.Item(8,3).Value = "鹏华稳健鸿利一年持有期混合A"
$ws.Cells.Item(8,4).Value = "2.61"
$ws.Cells.Item(8,5).Value = "92.98"
$ws.Cells.Item(8,6).Value = "6.04"
$ws.Cells.Item(8,7).Value = "0.1576"
$ws.Cells.Item(8,8).Value = 4

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "013392"
$ws.Cells.Item(9,3).Value = "广发港股通优质增长混合C"
$ws.Cells.Item(9,4).Value = "2.44"
$ws.Cells.Item(9,5).Value = "90.07"
$ws.Cells.Item(9,6).Value = "5.56"
$ws.Cells.Item(9,7).Value = "0.1357"
$ws.Cells.Item(9,8).Value = 10

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "008134"
$ws.Cells.Item(10,3).Value = "鹏华优选价值股票"
$ws.Cells.Item(10,4).Value = "1.80"
$ws.Cells.Item(10,5).Value = "92.72"
$ws.Cells.Item(10,6).Value = "5.45"
$ws.Cells.Item(10,7).Value = "0.0981"
$ws.Cells.Item(10,8).Value = 4

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "011157"
$ws.Cells.Item(11,3).Value = "弘毅远方港股通智选领航混合A"
$ws.Cells.Item(11,4).Value = "2.55"
$ws.Cells.Item(11,5).Value = "88.63"
$ws.Cells.Item(11,6).Value = "3.24"
$ws.Cells.Item(11,7).Value = "0.0826"
$ws.Cells.Item(11,8).Value = 10

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "011574"
$ws.Cells.Item(12,3).Value = "鹏华领航一年持有期混合A"
$ws.Cells.Item(12,4).Value = "1.20"
$ws.Cells.Item(12,5).Value = "92.84"
$ws.Cells.Item(12,6).Value = "5.55"
$ws.Cells.Item(12,7).Value = "0.0666"
$ws.Cells.Item(12,8).Value = 4

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "011575"
$ws.Cells.Item(13,3).Value = "鹏华领航一年持有期混合C"
$ws.Cells.Item(13,4).Value = "0.91"
$ws.Cells.Item(13,5).Value = "92.84"
$ws.Cells.Item(13,6).Value = "5.55"
$ws.Cells.Item(13,7).Value = "0.0505"
$ws.Cells.Item(13,8).Value = 4

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "011158"
$ws.Cells.Item(14,3).Value = "弘毅远方港股通智选领航混合C"
$ws.Cells.Item(14,4).Value = "0.75"
$ws.Cells.Item(14,5).Value = "88.63"
$ws.Cells.Item(14,6).Value = "3.24"
$ws.Cells.Item(14,7).Value = "0.0243"
$ws.Cells.Item(14,8).Value = 10

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "015187"
$ws.Cells.Item(15,3).Value = "汇添富消费升级混合C"
$ws.Cells.Item(15,4).Value = "0.52"
$ws.Cells.Item(15,5).Value = "87.74"
$ws.Cells.Item(15,6).Value = "2.36"
$ws.Cells.Item(15,7).Value = "0.0123"
$ws.Cells.Item(15,8).Value = 10

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "012641"
$ws.Cells.Item(16,3).Value = "鹏华稳健鸿利一年持有期混合C"
$ws.Cells.Item(16,4).Value = "0.10"
$ws.Cells.Item(16,5).Value = "92.98"
$ws.Cells.Item(16,6).Value = "6.04"
$ws.Cells.Item(16,7).Value = "0.0060"
$ws.Cells.Item(16,8).Value = 4

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "015183"
$ws.Cells.Item(17,3).Value = "汇添富文体娱乐主题混合C"
$ws.Cells.Item(17,4).Value = "0.01"
$ws.Cells.Item(17,5).Value = "90.70"
$ws.Cells.Item(17,6).Value = "2.92"
$ws.Cells.Item(17,7).Value = "0.0003"
$ws.Cells.Item(17,8).Value = 10

$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "015184"
$ws.Cells.Item(18,3).Value = "汇添富文体娱乐主题混合D"
$ws.Cells.Item(18,4).Value = "0.00"
$ws.Cells.Item(18,5).Value = "90.70"
$ws.Cells.Item(18,6).Value = "2.92"
# Last row's market value is a genuine number (0), not text like the rows above.
$ws.Cells.Item(18,7).NumberFormat = "General"
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = 10

# ---------------------------------------------------------------------------
# 2) Insert a new leading data row in "总计" for the 2022-Q3 figures
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Re-apply the bold/border/centered style used by the other index cells (A column)
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 17
$wsTotal.Cells.Item(2,4).Value = 7.07

Write-Host "2022-Q3 sheet added and 总计 updated"
